$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Another"): change the current selection ---
$ws2 = $wb.Worksheets.Item("Another")
[void]$ws2.Activate()
[void]$ws2.Range("A2:D4").Select()

# --- Add a new worksheet "NoHeaders" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "NoHeaders"

# Populate it with the same data as "Another", minus the header row.
# (Apply the date format before writing values so the existing date style
# is reused instead of Excel auto-creating a new one from the literal.)
$ws3.Range("D1:D3").NumberFormat = "d-mmm"

$ws3.Range("A1").Value = "a"
$ws3.Range("B1").Value = 1
$ws3.Range("C1").Value = $true
$ws3.Range("D1").Value = [DateTime]"2022-06-03"

$ws3.Range("A2").Value = "b"
$ws3.Range("B2").Value = 2
$ws3.Range("C2").Value = $false
$ws3.Range("D2").Value = [DateTime]"2022-05-12"

$ws3.Range("A3").Value = "c"
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = $false
$ws3.Range("D3").Value = [DateTime]"2022-02-15"

# Make it the active sheet/tab with a full-range selection.
[void]$ws3.Activate()
[void]$ws3.Range("A1:D3").Select()
